$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-99 down to 63-100
$ws.Rows.Item(62).Insert()

# Copy formatting/layout from the row now below (old row 62, now at 63) into the new blank row 62
$ws.Range("A63:R63").Copy()
$ws.Range("A62:R62").PasteSpecial()

# Populate the new row 62 with this week's price entry
$ws.Range("D62").Value = 45141
$ws.Range("J62").Value = 240
$ws.Range("K62").Value = 19000
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = 19500
$ws.Range("P62").Value = 780
